$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Account Information")

# New row 4: test / test / test@gmail.com / test / 12/12/1212 / 69176879
$ws.Cells.Item(4, 1).Value = "test"
$ws.Cells.Item(4, 2).Value = "test"
$ws.Cells.Item(4, 3).Value = "test@gmail.com"
$ws.Cells.Item(4, 4).Value = "test"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "12/12/1212"
$ws.Cells.Item(4, 5).ClearFormats()
$ws.Cells.Item(4, 6).Value = 69176879

# New row 5: notadmin / notadmin / it workssss / jack / 12-23-2222 / 87999924
$ws.Cells.Item(5, 1).Value = "notadmin"
$ws.Cells.Item(5, 2).Value = "notadmin"
$ws.Cells.Item(5, 3).Value = "it workssss"
$ws.Cells.Item(5, 4).Value = "jack"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "12-23-2222"
$ws.Cells.Item(5, 5).ClearFormats()
$ws.Cells.Item(5, 6).Value = 87999924
